$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Billing")
$ws.Activate()
$ws.Range("B2").Value = "Anila"
$ws.Range("C2").Value = "Niles"
$ws.Range("C2").Select()
